# Update crypto Price (D) and Volume(1h) (E) columns with latest values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.303.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "'1.873.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'0.7121"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'242.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.3114"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("D9").Value = "'0.07748"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").Value = "'25.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").Value = "'0.08474"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.85%  "
$ws.Range("D12").Value = "'1.874.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "'5.209"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'0.7120"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "'91.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "'0.000008382"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.68%  "
$ws.Range("D17").Value = "'29.306.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "'5.987"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "'242.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'13.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "'2.124.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'7.793"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'0.1619"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("D26").Value = "'162.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").Value = "'9.019"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").Value = "'1.508"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'4.416"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").Value = "'4.328"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.01%  "
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("D33").Value = "'0.05257"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("D34").Value = "'1.921"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'1.173"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'0.7430"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("D37").Value = "'2.682"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'0.01859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "'2.718"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").Value = "'1.165.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").Value = "'6.345"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.42%  "
$ws.Range("D42").Value = "'0.8892"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").Value = "'106.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.00%  "
$ws.Range("D45").Value = "'0.9994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'2.019.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "'1.808"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.56%  "
$ws.Range("D48").Value = "'0.5196"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").Value = "'9.387"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("D51").Value = "'0.4300"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.58%  "
